$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 9 (August) label and values
$ws.Range("A9").Value = "August (through 08-24)"
$ws.Range("B9").Value = 26
$ws.Range("C9").Value = 55
$ws.Range("D9").Value = 68
$ws.Range("E9").Value = 43
$ws.Range("F9").Value = 35
$ws.Range("G9").Value = 138
$ws.Range("H9").Value = 121

# Update row 10 (Total) values
$ws.Range("B10").Value = 188
$ws.Range("C10").Value = 357
$ws.Range("D10").Value = 533
$ws.Range("E10").Value = 468
$ws.Range("F10").Value = 339
$ws.Range("G10").Value = 759
$ws.Range("H10").Value = 1036
